$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.013.29"
$ws.Range("E2").Value = "  +0.12%  "

$ws.Range("D3").Value = "3.525.07"
$ws.Range("E3").Value = "  -0.79%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.16"
$ws.Range("E5").Value = "  -1.05%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.37"
$ws.Range("E6").Value = "  -1.73%  "

$ws.Range("D7").Value = "3.523.12"
$ws.Range("E7").Value = "  -0.80%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.492"
$ws.Range("E9").Value = "  -1.03%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.125"
$ws.Range("E10").Value = "  +0.98%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.15"
$ws.Range("E11").Value = "  +3.27%  "

$ws.Range("E12").Value = "  +0.30%  "

$ws.Range("D13").Value = "4.126.49"
$ws.Range("E13").Value = "  -0.67%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.68"
$ws.Range("E14").Value = "  +1.96%  "

$ws.Range("E15").Value = "  -0.16%  "

$ws.Range("E16").Value = "  +0.50%  "

$ws.Range("D17").Value = "3.523.84"
$ws.Range("E17").Value = "  -0.88%  "

$ws.Range("D18").Value = "65.032.70"
$ws.Range("E18").Value = "  +0.27%  "

$ws.Range("E19").Value = "  +1.04%  "

$ws.Range("E20").Value = "  +0.44%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.73"
$ws.Range("E21").Value = "  -1.97%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "392.50"
$ws.Range("E22").Value = "  +0.83%  "

$ws.Range("E23").Value = "  +0.70%  "

$ws.Range("D24").Value = "3.667.38"
$ws.Range("E24").Value = "  -0.86%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.77"
$ws.Range("E25").Value = "  +0.82%  "

$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("E27").Value = "  -3.15%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.75"
$ws.Range("E28").Value = "  +1.21%  "

$ws.Range("E29").Value = "  +10.64%  "

$ws.Range("E30").Value = "  -0.06%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.28"
$ws.Range("E31").Value = "  -0.63%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.38"
$ws.Range("E32").Value = "  +0.49%  "

$ws.Range("D33").Value = "3.529.20"
$ws.Range("E33").Value = "  -0.95%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "24.15"
$ws.Range("E34").Value = "  +0.56%  "

$ws.Range("E35").Value = "  +0.00%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.145"
$ws.Range("E36").Value = "  -0.67%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.34"
$ws.Range("E37").Value = "  +6.81%  "

$ws.Range("E38").Value = "  +3.07%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.03"
$ws.Range("E39").Value = "  +1.19%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "168.72"
$ws.Range("E40").Value = "  -0.94%  "

$ws.Range("E41").Value = "  +1.09%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.825"
$ws.Range("E42").Value = "  -0.34%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.28"
$ws.Range("E43").Value = "  +5.29%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "25.91"
$ws.Range("E44").Value = "  -2.88%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.97"
$ws.Range("E45").Value = "  +0.78%  "

$ws.Range("E46").Value = "  +0.04%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.45"
$ws.Range("E47").Value = "  -0.14%  "

$ws.Range("E48").Value = "  +0.59%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.95"
$ws.Range("E49").Value = "  +0.82%  "

$ws.Range("D50").Value = "2.425.56"
$ws.Range("E50").Value = "  -0.86%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.911"
$ws.Range("E51").Value = "  +6.35%  "
